$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the "Make GitHub repository" task row (old row 7).
$ws.Rows("7:7").Delete()

# 2) Row 8 ("Fetch and put the IPs in our new Database") gains an hours note
#    ("??") in column C and a comment in column D.
$ws.Range("C8").Value = "??"
$ws.Range("D8").Value = "Isn't it the automated thing we do everyday? So basically no time from us, but the scripts are working instead"

# 3) Fill in newly logged hours for a task above the insertion point below.
$ws.Range("C17").Value = 1

# 4) Insert a new task row after "Configure SSH on Shreyas' laptop" (row 18)
#    for the new "Create scrpt to fetch blocked Ips" task.
$ws.Rows("19:19").Insert()
$ws.Range("A19").Value = "Create scrpt to fetch blocked Ips"
$ws.Range("B19").Value = "Márton Reiter"
$ws.Range("C19").Value = 5

# Match the copied row's wrap-text style used by the other task-name/participant cells.
$ws.Range("A19").Style = $ws.Range("A18").Style
$ws.Range("B19").Style = $ws.Range("B18").Style
$ws.Range("C19").Style = $ws.Range("C17").Style

# 5) Fill in newly logged hours for tasks below the insertion point (rows
#    renumbered after the insert above).
$ws.Range("C41").Value = 5
$ws.Range("C48").Value = 3

# Keep the selection where the author left it.
$ws.Range("D17").Select()
